# "Add files via upload" - refreshed data values for Tab09 (Inequality/Poverty
# indicators) aggregate rows: Reste du monde, Amerique latine et Caraibes,
# Asie (hors revenu eleve), Monde, ASEAN, Afrique a faible revenu, RDM a
# revenu intermediaire superieur, Pays a revenu eleve.
# Only the numeric columns C:G (poverty-rate / Gini / income-share
# aggregates) were updated; everything else (labels, formatting, styles)
# is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 63 - Reste du monde
$ws.Range("C63").Value = 2.1649484536082499
$ws.Range("D63").Value = 7.8154639175257801
$ws.Range("E63").Value = 21.2144329896907
$ws.Range("F63").Value = 34.550515463917499
$ws.Range("G63").Value = 42.231958762886599

# Row 64 - Amerique latine et Caraibes
$ws.Range("C64").Value = 2.6380952380952398
$ws.Range("D64").Value = 7.4666666666666703
$ws.Range("F64").Value = 44.185714285714297
$ws.Range("G64").Value = 49.985714285714302

# Row 65 - Asie (pays a revenu eleve exclus)
$ws.Range("C65").Value = 5.35
$ws.Range("D65").Value = 20.936363636363598
$ws.Range("E65").Value = 49.072727272727299
$ws.Range("F65").Value = 33.6727272727273
$ws.Range("G65").Value = 41.909090909090899

# Row 66 - Monde
$ws.Range("C66").Value = 10.8090909090909
$ws.Range("D66").Value = 22.371328671328701
$ws.Range("E66").Value = 39.316083916083898
$ws.Range("F66").Value = 36.558041958041997
$ws.Range("G66").Value = 43.956643356643397

# Row 76 - ASEAN
$ws.Range("C76").Value = 2.1428571428571401
$ws.Range("D76").Value = 13.271428571428601
$ws.Range("E76").Value = 41
$ws.Range("F76").Value = 37.171428571428599
$ws.Range("G76").Value = 44.771428571428601

# Row 83 - RDM (pays riches en ressources exclus)
$ws.Range("C83").Value = 1.8218390804597699
$ws.Range("D83").Value = 7
$ws.Range("E83").Value = 20.160919540229902
$ws.Range("F83").Value = 34.628735632183897
$ws.Range("G83").Value = 42.278160919540198

# Row 89 - RDM, pays a revenu intermediaire, tranche superieure
$ws.Range("C89").Value = 1.49714285714286
$ws.Range("D89").Value = 5.5457142857142898
$ws.Range("E89").Value = 21.5857142857143
$ws.Range("F89").Value = 37.851428571428599
$ws.Range("G89").Value = 44.9514285714286

# Row 90 - Pays a revenu eleve
$ws.Range("C90").Value = 0.30249999999999999
$ws.Range("D90").Value = 0.61499999999999999
$ws.Range("F90").Value = 31.93
$ws.Range("G90").Value = 39.880000000000003
